# Update the "Correspond Handoff Datetime" / "Correspond Handback DateTime"
# timestamps for the 879ee5fc-... row (row 3) on the zh-cn and de-de
# report sheets, to reflect a re-generated handback report.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E3").Value = "2016-03-24 00:47:10"
$wsZhCn.Range("H3").Value = "2016-03-24 00:47:37"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E3").Value = "2016-03-24 00:47:14"
$wsDeDe.Range("H3").Value = "2016-03-24 00:47:43"
